# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect freshly generated output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# 展览 sheet updates
$wsExhibit.Range("F2").Value  = 11408
$wsExhibit.Range("F3").Value  = 10768
$wsExhibit.Range("F6").Value  = 979
$wsExhibit.Range("F9").Value  = 36
$wsExhibit.Range("F11").Value = 10571
$wsExhibit.Range("F17").Value = 102
$wsExhibit.Range("F19").Value = 11091
$wsExhibit.Range("F20").Value = 10845

# 全部类型 sheet updates
$wsAll.Range("F2").Value  = 11408
$wsAll.Range("F3").Value  = 10769
$wsAll.Range("F6").Value  = 979
$wsAll.Range("F9").Value  = 36
$wsAll.Range("F11").Value = 0
$wsAll.Range("F17").Value = 102
$wsAll.Range("F19").Value = 11091
$wsAll.Range("F20").Value = 10845
